$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (Table1 row): PERT estimates a/m/b updated -> a=3, m=4, b=5 (TE recalculates to 4)
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 4
$ws.Range("G25").Value = 5

# Row 29 (Table1 row): PERT estimates m/b updated -> m=2, b=2 (TE recalculates to 2)
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 2

# Scroll the viewport and move the active selection to match the saved view.
$ws.Range("G25").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 2
